$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets   = $wb.Worksheets.Item("Assets")

# --- Settings sheet: fill in rows 6-9 with the robot-registration asset rows ---
$wsSettings.Range("A6").Value = "GoogleFormLink"
$wsSettings.Range("B6").Value = "GoogleFormLink"
$wsSettings.Range("C6").Value = "This is link to google form that any user of trip-planner need to fulfill to provide"

$wsSettings.Range("A7").Value = "LetterSubject"
$wsSettings.Range("B7").Value = "LetterSubject"
$wsSettings.Range("C7").Value = "Trip-Planner Robot Registration Response Letter Subject "

$wsSettings.Range("A8").Value = "LetterText"
$wsSettings.Range("B8").Value = "LetterText"
$wsSettings.Range("C8").Value = "Text that will be send for each user who would like to start registration in trip-planner"

$wsSettings.Range("A9").Value = "EmailCredentials"
$wsSettings.Range("B9").Value = "EmailCredentials"
$wsSettings.Range("C9").Value = "This is credentials for email in which we receive registration letters from users, and from which we send responses and new data about trips."

# --- Assets sheet: clear out the old values now that they live on Settings ---
$wsAssets.Range("A2:B2").ClearContents()
$wsAssets.Range("C2:D2").ClearContents()

$wsAssets.Range("A3:B3").ClearContents()
$wsAssets.Range("C3").ClearContents()

$wsAssets.Range("A4:B4").ClearContents()
$wsAssets.Range("C4:D4").ClearContents()

$wsAssets.Range("A5:D5").ClearContents()

# --- Selections, as left by the editor ---
$wsSettings.Range("B18").Select()
$wsAssets.Range("C9").Select()
